$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.860.23"
$ws.Range("E2").Value = "  +2.07%  "

$ws.Range("D3").Value = "3.589.30"
$ws.Range("E3").Value = "  +1.33%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.43"
$ws.Range("E5").Value = "  +3.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.94"
$ws.Range("E6").Value = "  +1.13%  "

$ws.Range("D7").Value = "3.585.30"
$ws.Range("E7").Value = "  +1.43%  "

$ws.Range("E8").Value = "  +1.07%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.200"
$ws.Range("E10").Value = "  +4.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.47"
$ws.Range("E11").Value = "  +9.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.591"
$ws.Range("E12").Value = "  +0.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.32"
$ws.Range("E13").Value = "  -0.46%  "

$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("D15").Value = "4.170.52"
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.46"
$ws.Range("E16").Value = "  -1.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "620.80"
$ws.Range("E17").Value = "  -1.50%  "

$ws.Range("D18").Value = "3.580.10"
$ws.Range("E18").Value = "  +1.08%  "

$ws.Range("D19").Value = "71.042.94"
$ws.Range("E19").Value = "  +2.27%  "

$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.52"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.37"
$ws.Range("E23").Value = "  -16.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.14"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.09"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.66"
$ws.Range("E28").Value = "  +0.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.37"
$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("E30").Value = "  +2.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.55"
$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.08"
$ws.Range("E32").Value = "  -2.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.18"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.31"
$ws.Range("E34").Value = "  -2.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "634.06"
$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.81"
$ws.Range("E36").Value = "  +8.16%  "

$ws.Range("E37").Value = "  -0.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.90"
$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0486"
$ws.Range("E39").Value = "  +6.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "57.41"
$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("E42").Value = "  +4.26%  "

$ws.Range("D43").Value = "3.421.17"
$ws.Range("E43").Value = "  +0.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.326"
$ws.Range("E44").Value = "  -1.28%  "

$ws.Range("E45").Value = "  +2.42%  "

$ws.Range("E47").Value = "  +5.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.14"
$ws.Range("E48").Value = "  +0.71%  "

$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.96"
$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("E51").Value = "  -0.03%  "
